# chore(runtime): publish files + archive (2025-12-03 18:09:51)
#
# Adds 3 newly played KHL matches to Matches_SOG, and refreshes the
# as_of_utc timestamps + rolling shots-on-goal aggregates on Shots_HA /
# Shots_Summary / Meta_ext to reflect them.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Matches_SOG: append rows 345-347 (3 new finished matches)
# ---------------------------------------------------------------------
$matches = $wb.Worksheets.Item("Matches_SOG")

$newRows = @(
    @{ Row = 345; Uid = "897838"; Date = "2025-12-02T12:15:00"; Home = "Амур";      Away = "Динамо Мн"; Sog = 27; Sga = 43 },
    @{ Row = 346; Uid = "897839"; Date = "2025-12-02T12:30:00"; Home = "Адмирал";   Away = "ХК Сочи";   Sog = 35; Sga = 23 },
    @{ Row = 347; Uid = "897840"; Date = "2025-12-02T19:00:00"; Home = "Локомотив"; Away = "СКА";       Sog = 48; Sga = 29 }
)

foreach ($m in $newRows) {
    $r = $m.Row

    # uid looks numeric ("897838") but must stay text like every other
    # row in this column -- force text format, write, then strip the
    # number-format back off so no stray style is left on the cell.
    $cellA = $matches.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $m.Uid
    $cellA.ClearFormats()

    $matches.Cells.Item($r, 2).Value = $m.Date
    $matches.Cells.Item($r, 3).Value = $m.Home
    $matches.Cells.Item($r, 4).Value = $m.Away
    $matches.Cells.Item($r, 5).Value = $m.Sog
    $matches.Cells.Item($r, 6).Value = $m.Sga
    $matches.Cells.Item($r, 7).Value = "khl_text"
}

# ---------------------------------------------------------------------
# Shots_HA: bump as_of_utc for every team, then refresh the home/away
# shot totals + per-game rates for the 6 teams involved in tonight's
# matches (Амур/Динамо Мн, Адмирал/ХК Сочи, Локомотив/СКА).
# ---------------------------------------------------------------------
$shotsHA = $wb.Worksheets.Item("Shots_HA")
$shotsHA.Range("D2:D23").Value = "2025-12-02T19:00:00Z"

# row => { GP_home, GP_away, HOGF_total, HOGA_total, HOGF_pg, HOGA_pg, AOGF_total, AOGA_total, AOGF_pg, AOGA_pg }
$haUpdates = @{
    4  = @{ E = 14; F = 16; G = 522; H = 375; I = 37.3; J = 26.8 }                                  # Адмирал (home)
    6  = @{ E = 16; F = 16; G = 485; H = 554; I = 30.3; J = 34.6 }                                  # Амур (home)
    9  = @{ F = 14; K = 493; L = 393; M = 35.2; N = 28.1 }                                          # Динамо Мн (away)
    12 = @{ E = 17; F = 17; G = 572; H = 450; I = 33.6; J = 26.5 }                                  # Локомотив (home)
    15 = @{ F = 15; K = 443; L = 508; M = 29.5; N = 33.9 }                                          # СКА (away)
    22 = @{ F = 15; K = 392; L = 557; M = 26.1; N = 37.1 }                                          # ХК Сочи (away)
}

foreach ($row in $haUpdates.Keys) {
    $vals = $haUpdates[$row]
    foreach ($col in $vals.Keys) {
        $shotsHA.Range("$col$row").Value = $vals[$col]
    }
}

# ---------------------------------------------------------------------
# Shots_Summary: same as_of_utc bump, then refresh combined SOG totals
# + per-game rates for the same 6 teams.
# ---------------------------------------------------------------------
$shotsSummary = $wb.Worksheets.Item("Shots_Summary")
$shotsSummary.Range("D2:D23").Value = "2025-12-02T19:00:00Z"

# row => { GP_total, SOG_total, SOGA_total, SOG_pg, SOGA_pg }
$summaryUpdates = @{
    4  = @{ E = 30; F = 1019; G = 819;  H = 34;   I = 27.3 }   # Адмирал
    6  = @{ E = 32; F = 923;  G = 1168; H = 28.8; I = 36.5 }   # Амур
    9  = @{ E = 32; F = 1154; G = 876;  H = 36.1 }             # Динамо Мн
    12 = @{ E = 34; F = 1096; G = 870;  H = 32.2; I = 25.6 }   # Локомотив
    15 = @{ E = 30; F = 931;  G = 1004; H = 31;   I = 33.5 }   # СКА
    22 = @{ E = 31; F = 864;  G = 1069; H = 27.9 }             # ХК Сочи
}

foreach ($row in $summaryUpdates.Keys) {
    $vals = $summaryUpdates[$row]
    foreach ($col in $vals.Keys) {
        $shotsSummary.Range("$col$row").Value = $vals[$col]
    }
}

# ---------------------------------------------------------------------
# Meta_ext: bump as_of_utc + build_version
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Meta_ext")
$meta.Range("B2").Value = "2025-12-02T19:00:00Z"
$meta.Range("D2").Value = 29
